$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the TOR450 rows (23-27) first so row numbers for TOR130 stay valid,
# then delete the TOR130 rows (2-8). This leaves only the TOR330 data
# (originally rows 9-22), which will shift up to rows 2-15.
$ws.Range("A23:D27").EntireRow.Delete() | Out-Null
$ws.Range("A2:D8").EntireRow.Delete() | Out-Null
